# Apply the timesheet update:
#  - Resource name:        Vineet Rajput            -> Tanuj Khaturia
#  - Week range:            7/23/2018-7/29/2018      -> 7/30/2018-8/5/2018
#  - Daily date headers:    23-Jul..29-Jul            -> 30-Jul,31-Jul,1-Aug..5-Aug
#  - Activity description: Client Call ( sanchit )   -> Client Call ( 1 )

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Tanuj Khaturia"

$ws.Range("B5").Value = "7/30/2018-8/5/2018"

$ws.Range("E7").Value = "30-Jul"
$ws.Range("F7").Value = "31-Jul"
$ws.Range("G7").Value = "1-Aug"
$ws.Range("H7").Value = "2-Aug"
$ws.Range("I7").Value = "3-Aug"
$ws.Range("J7").Value = "4-Aug"
$ws.Range("K7").Value = "5-Aug"

$ws.Range("B9").Value = "Client Call ( 1 ) "
